$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C26").Value = "Dabney things"
$ws.Range("C27").Value = "Lewis Carol"
$ws.Range("C28").Value = "Art of computing"

$ws.Range("C29").Select()
